$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "63.344.51"
$ws.Range("E2").Value = "  +3.96%  "
$ws.Range("D3").Value = "3.485.07"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("D5").Value = "584.27"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").Value = "147.74"
$ws.Range("E6").Value = "  +6.71%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("D9").Value = "7.71"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("D12").Value = "4.081.27"
$ws.Range("E12").Value = "  +3.45%  "
$ws.Range("D13").Value = "29.72"
$ws.Range("E13").Value = "  +5.95%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "3.483.78"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "63.341.17"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("D18").Value = "6.29"
$ws.Range("E18").Value = "  +3.28%  "
$ws.Range("D19").Value = "14.37"
$ws.Range("E19").Value = "  +6.34%  "
$ws.Range("D20").Value = "9.36"
$ws.Range("E20").Value = "  +5.09%  "
$ws.Range("D21").Value = "391.34"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  +7.36%  "
$ws.Range("D26").Value = "3.627.50"
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("E28").Value = "  +9.94%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "8.29"
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  +7.50%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "23.79"
$ws.Range("E34").Value = "  +3.40%  "
$ws.Range("D35").Value = "32.67"
$ws.Range("E35").Value = "  +27.99%  "
$ws.Range("D36").Value = "5.34"
$ws.Range("E36").Value = "  +8.71%  "
$ws.Range("D37").Value = "7.14"
$ws.Range("E37").Value = "  +4.59%  "
$ws.Range("D38").Value = "171.78"
$ws.Range("E39").Value = "  +8.80%  "
$ws.Range("D40").Value = "3.522.25"
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("D41").Value = "0.0768"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E42").Value = "  +4.63%  "
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "1.73"
$ws.Range("E44").Value = "  +6.41%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "42.46"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  +9.32%  "
$ws.Range("D47").Value = "2.628.64"
$ws.Range("E47").Value = "  +7.49%  "
$ws.Range("D48").Value = "23.69"
$ws.Range("E48").Value = "  +7.37%  "
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").Value = "  +15.33%  "
$ws.Range("D50").Value = "6.76"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  +5.02%  "
